# Matriz de trazabilidad - actualizacion de fechas de estado, version y
# estado de avance (Visualizar Resenas + Gestionar Menu + Actualizacion
# Matriz de Trazabilidad).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# ---------------------------------------------------------------------
# Columna I (Fecha de estado): la mayoria de las filas pasan del
# 23/10/2020 (44127) al 25/10/2020 (44129); la fila 27 pasa al 24/10/2020
# (44128).
# ---------------------------------------------------------------------
$rowsTo44129 = @(3,4,5,6,7,8,9,10,11,12,13,14,15,20,21,26,28,30,31,32,33,34,35,36,37,38,39)
foreach ($r in $rowsTo44129) {
    $ws.Range("I$r").Value2 = 44129
}
$ws.Range("I27").Value2 = 44128

# ---------------------------------------------------------------------
# Columna G (Version): RF-09 (fila 13) pasa de v1 a v2, mientras que
# RF-10 (fila 28) pasa de v2 a v1.
# ---------------------------------------------------------------------
$ws.Range("G13").Value = "v2"
$ws.Range("G28").Value = "v1"

# ---------------------------------------------------------------------
# Columna H (Estado): filas 26 y 27 pasan de "Falta" a "Completado".
# Se copia el formato (relleno/estilo "Bueno") de una celda que ya tiene
# el estado "Completado" para mantener el mismo estilo condicional.
# ---------------------------------------------------------------------
$ws.Range("H16").Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4122) | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("H26").Value = "Completado"
$ws.Range("H27").Value = "Completado"

# ---------------------------------------------------------------------
# Estado de la ventana / seleccion: la hoja quedo desplazada hacia la
# fila 31 y con la celda R13 seleccionada.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
$win.Zoom = 60
$ws.Range("R13").Select()
